$wb = $excel.ActiveWorkbook

# ----- FIRE sheet (index 1) -----
$ws1 = $wb.Worksheets.Item(1)

# Rows 2-7: F = 1 (G,H unchanged)
$ws1.Range("F2:F7").Value = 1

# Rows 8-10: F=1, G=0.08, H=0.32
$ws1.Range("F8:F10").Value = 1
$ws1.Range("G8:G10").Value = 0.08
$ws1.Range("H8:H10").Value = 0.32

# Rows 11-13: F=1, G=0.15, H=0
$ws1.Range("F11:F13").Value = 1
$ws1.Range("G11:G13").Value = 0.15
$ws1.Range("H11:H13").Value = 0

# Rows 14-17: F=1, G=0.06, H=0.48
$ws1.Range("F14:F17").Value = 1
$ws1.Range("G14:G17").Value = 0.06
$ws1.Range("H14:H17").Value = 0.48

# Rows 18-20: F = 1 (G,H unchanged)
$ws1.Range("F18:F20").Value = 1

# Rows 21-23: F=1.2, G=0.075, H=0.6
$ws1.Range("F21:F23").Value = 1.2
$ws1.Range("G21:G23").Value = 0.075
$ws1.Range("H21:H23").Value = 0.6

# ----- WATER sheet (index 2) -----
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("F2:F23").Value = 1

# ----- NATURE sheet (index 3) -----
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("F2:F25").Value = 1

# ----- View state: selection / active sheet -----
# WATER selection: F2:F23, active cell F2 (visited first so it is not left as tabSelected)
$ws2.Activate() | Out-Null
$ws2.Range("F2:F23").Select() | Out-Null

# NATURE selection: F2:F25, active cell F2
$ws3.Activate() | Out-Null
$ws3.Range("F2:F25").Select() | Out-Null

# FIRE is the final active/selected sheet, active cell F25
$ws1.Activate() | Out-Null
$ws1.Range("F25").Select() | Out-Null
